# Commit message: "adding another graphic in the last screen"
# The author cleared out the sample/demo data rows (rows 2-5) that were
# previously used to preview the chart, leaving just the header row and the
# formatted trailing row (row 14) intact, then left the selection on A2
# (presumably right before inserting another chart/graphic).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of the now-unneeded sample rows (2 through 5),
# keeping row 1 (headers) and row 14 untouched.
$ws.Range("A2:L5").ClearContents()

# Move/leave the active selection at A2, matching the saved view state.
$ws.Range("A2").Select()
